# Updated symbol list (Price / Volume(1h) columns) to match refreshed
# crypto price feed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes $newValue into $cellRef as plain text (no auto number/percent
# conversion), preserving the original "General" style (no style index)
# on the cell once done.
function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}

Set-TextValue "D2" "278.77"
Set-TextValue "E2" "1.02%"
Set-TextValue "D3" "27.37"
Set-TextValue "E3" "0.02%"
Set-TextValue "D4" "4.826"
Set-TextValue "E4" "0.39%"
Set-TextValue "E5" "0.40%"
Set-TextValue "D6" "7.033"
Set-TextValue "E6" "1.16%"
Set-TextValue "E7" "-2.51%"
Set-TextValue "D8" "0.8916"
Set-TextValue "E8" "1.41%"
Set-TextValue "D9" "0.1526"
Set-TextValue "E9" "0.24%"
Set-TextValue "D10" "0.05632"
Set-TextValue "E10" "11.28%"
Set-TextValue "D11" "0.07499"
Set-TextValue "E11" "-0.19%"
Set-TextValue "D12" "0.02913"
Set-TextValue "E12" "-2.19%"
Set-TextValue "D13" "0.08979"
Set-TextValue "E13" "-0.48%"
Set-TextValue "D14" "0.001576"
Set-TextValue "E14" "0.87%"
Set-TextValue "D15" "0.0006386"
Set-TextValue "E15" "-0.41%"
Set-TextValue "D16" "0.006094"
Set-TextValue "E16" "4.78%"
Set-TextValue "D17" "3.472"
Set-TextValue "E17" "0.63%"
Set-TextValue "D18" "3.302"
Set-TextValue "E18" "-0.07%"
Set-TextValue "D19" "2.288"
Set-TextValue "E19" "0.16%"
Set-TextValue "D20" "0.3147"
Set-TextValue "E20" "0.92%"
Set-TextValue "E21" "-0.31%"
Set-TextValue "D22" "3.905"
Set-TextValue "E22" "-1.42%"
Set-TextValue "D23" "0.1505"
Set-TextValue "E23" "9.04%"
Set-TextValue "D24" "0.04380"
Set-TextValue "E24" "-0.81%"
Set-TextValue "D25" "0.001176"
Set-TextValue "E25" "0.23%"
Set-TextValue "D26" "0.004281"
Set-TextValue "E26" "10.77%"
Set-TextValue "D28" "0.0001177"
Set-TextValue "E28" "-1.83%"
Set-TextValue "E29" "-14.72%"
Set-TextValue "D40" "0.04038"
Set-TextValue "E40" "-3.25%"
Set-TextValue "D41" "0.006719"
Set-TextValue "E41" "-1.98%"
Set-TextValue "D42" "0.1405"
Set-TextValue "E42" "19.44%"
Set-TextValue "D43" "0.002046"
Set-TextValue "E43" "0.82%"
Set-TextValue "E44" "-3.26%"
Set-TextValue "D45" "0.00005537"
Set-TextValue "E45" "6.94%"
Set-TextValue "D46" "1.561"
Set-TextValue "E46" "5.00%"
Set-TextValue "D47" "0.01846"
Set-TextValue "E47" "-19.69%"
